# Update "想去人数" (column F) values across sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 16   # F3: 15 -> 16
$ws.Cells.Item(5, 6).Value = 4672   # F5: 4661 -> 4672
$ws.Cells.Item(7, 6).Value = 128   # F7: 127 -> 128
$ws.Cells.Item(8, 6).Value = 105   # F8: 0 -> 105
$ws.Cells.Item(9, 6).Value = 95   # F9: 94 -> 95
$ws.Cells.Item(11, 6).Value = 0   # F11: 714 -> 0
$ws.Cells.Item(12, 6).Value = 198   # F12: 0 -> 198
$ws.Cells.Item(13, 6).Value = 1060   # F13: 1057 -> 1060
$ws.Cells.Item(14, 6).Value = 89   # F14: 0 -> 89
$ws.Cells.Item(15, 6).Value = 253   # F15: 252 -> 253
$ws.Cells.Item(17, 6).Value = 75   # F17: 0 -> 75
$ws.Cells.Item(18, 6).Value = 129   # F18: 127 -> 129
$ws.Cells.Item(19, 6).Value = 0   # F19: 100 -> 0
$ws.Cells.Item(20, 6).Value = 0   # F20: 3661 -> 0
$ws.Cells.Item(22, 6).Value = 0   # F22: 37 -> 0
$ws.Cells.Item(24, 6).Value = 84   # F24: 0 -> 84
$ws.Cells.Item(33, 6).Value = 129   # F33: 127 -> 129
$ws.Cells.Item(34, 6).Value = 240   # F34: 236 -> 240
$ws.Cells.Item(35, 6).Value = 281   # F35: 279 -> 281
$ws.Cells.Item(37, 6).Value = 141   # F37: 0 -> 141
$ws.Cells.Item(38, 6).Value = 1536   # F38: 1535 -> 1536
$ws.Cells.Item(39, 6).Value = 924   # F39: 923 -> 924
$ws.Cells.Item(40, 6).Value = 27   # F40: 26 -> 27
$ws.Cells.Item(41, 6).Value = 35   # F41: 34 -> 35
$ws.Cells.Item(42, 6).Value = 53   # F42: 0 -> 53
$ws.Cells.Item(43, 6).Value = 465   # F43: 0 -> 465
$ws.Cells.Item(44, 6).Value = 0   # F44: 475 -> 0
$ws.Cells.Item(45, 6).Value = 0   # F45: 70 -> 0
# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 0   # F2: 100 -> 0
# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 0   # F2: 46 -> 0
$ws.Cells.Item(4, 6).Value = 203   # F4: 200 -> 203
$ws.Cells.Item(5, 6).Value = 4672   # F5: 4661 -> 4672
$ws.Cells.Item(7, 6).Value = 0   # F7: 127 -> 0
$ws.Cells.Item(8, 6).Value = 105   # F8: 0 -> 105
$ws.Cells.Item(10, 6).Value = 95   # F10: 0 -> 95
$ws.Cells.Item(12, 6).Value = 0   # F12: 714 -> 0
$ws.Cells.Item(13, 6).Value = 0   # F13: 197 -> 0
$ws.Cells.Item(14, 6).Value = 1060   # F14: 1057 -> 1060
$ws.Cells.Item(15, 6).Value = 0   # F15: 89 -> 0
$ws.Cells.Item(16, 6).Value = 253   # F16: 252 -> 253
$ws.Cells.Item(17, 6).Value = 161   # F17: 0 -> 161
$ws.Cells.Item(18, 6).Value = 0   # F18: 74 -> 0
$ws.Cells.Item(19, 6).Value = 129   # F19: 127 -> 129
$ws.Cells.Item(20, 6).Value = 101   # F20: 100 -> 101
$ws.Cells.Item(21, 6).Value = 3674   # F21: 0 -> 3674
$ws.Cells.Item(22, 6).Value = 6016   # F22: 5994 -> 6016
$ws.Cells.Item(26, 6).Value = 528   # F26: 0 -> 528
$ws.Cells.Item(27, 6).Value = 0   # F27: 45 -> 0
$ws.Cells.Item(28, 6).Value = 3405   # F28: 0 -> 3405
$ws.Cells.Item(29, 6).Value = 377   # F29: 376 -> 377
$ws.Cells.Item(30, 6).Value = 30   # F30: 0 -> 30
$ws.Cells.Item(32, 6).Value = 0   # F32: 567 -> 0
$ws.Cells.Item(33, 6).Value = 0   # F33: 519 -> 0
$ws.Cells.Item(34, 6).Value = 129   # F34: 127 -> 129
$ws.Cells.Item(35, 6).Value = 240   # F35: 236 -> 240
$ws.Cells.Item(36, 6).Value = 0   # F36: 279 -> 0
$ws.Cells.Item(37, 6).Value = 360   # F37: 0 -> 360
$ws.Cells.Item(38, 6).Value = 0   # F38: 137 -> 0
$ws.Cells.Item(39, 6).Value = 1536   # F39: 1535 -> 1536
$ws.Cells.Item(40, 6).Value = 924   # F40: 923 -> 924
$ws.Cells.Item(41, 6).Value = 0   # F41: 26 -> 0
$ws.Cells.Item(42, 6).Value = 35   # F42: 34 -> 35
$ws.Cells.Item(43, 6).Value = 0   # F43: 51 -> 0
$ws.Cells.Item(44, 6).Value = 0   # F44: 465 -> 0
$ws.Cells.Item(47, 6).Value = 0   # F47: 556 -> 0
